$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

# Update attendance note for 9/29 (D9 changed from "will" to "Younouss")
$ws.Range("D9").Value = "Younouss"

# Updated status/notes for the backlog items
$ws.Range("F20").Value = "done"

$ws.Range("F22").Value = "in progress"
$ws.Range("G22").Value = "enemies do not currently attack"

$ws.Range("G23").Value = "pending on difficulty settings"

$ws.Range("F24").Value = "done"

$ws.Range("F25").Value = "in progress"
$ws.Range("G25").Value = "image file that is in use"

# Update selection / scroll position
$null = $ws.Range("D11").Select()
